# Weekly Fruta/Hortaliza update:
# Insert 3 new rows (Perú "$/bandeja 10 kilos" price points for a new
# reporting date) at the top of the data block, pushing the existing
# records down by 3 rows (891-925 -> 894-928).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(891).Insert()
$ws.Rows.Item(892).Insert()
$ws.Rows.Item(893).Insert()

# --- Row 891: Especial ---
$ws.Cells.Item(891, 1).Value = 11
$ws.Cells.Item(891, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(891, 3).Value = "Bíobío"
$ws.Cells.Item(891, 4).Value = 45041
$ws.Cells.Item(891, 5).Value = 8
$ws.Cells.Item(891, 6).Value = "Fruta"
$ws.Cells.Item(891, 7).Value = 100106
$ws.Cells.Item(891, 8).Value = "Oleaginosos"
$ws.Cells.Item(891, 9).Value = 100106002
$ws.Cells.Item(891, 10).Value = "Palta"
$ws.Cells.Item(891, 11).Value = "Hass"
$ws.Cells.Item(891, 12).Value = "Especial"
$ws.Cells.Item(891, 13).Value = 50
$ws.Cells.Item(891, 14).Value = 35000
$ws.Cells.Item(891, 15).Value = 35000
$ws.Cells.Item(891, 16).Value = 35000
$ws.Cells.Item(891, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(891, 18).Value = "Perú"
$ws.Cells.Item(891, 19).Value = 3500
$ws.Cells.Item(891, 20).Value = 10

# --- Row 892: Primera ---
$ws.Cells.Item(892, 1).Value = 11
$ws.Cells.Item(892, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(892, 3).Value = "Bíobío"
$ws.Cells.Item(892, 4).Value = 45041
$ws.Cells.Item(892, 5).Value = 8
$ws.Cells.Item(892, 6).Value = "Fruta"
$ws.Cells.Item(892, 7).Value = 100106
$ws.Cells.Item(892, 8).Value = "Oleaginosos"
$ws.Cells.Item(892, 9).Value = 100106002
$ws.Cells.Item(892, 10).Value = "Palta"
$ws.Cells.Item(892, 11).Value = "Hass"
$ws.Cells.Item(892, 12).Value = "Primera"
$ws.Cells.Item(892, 13).Value = 50
$ws.Cells.Item(892, 14).Value = 30000
$ws.Cells.Item(892, 15).Value = 30000
$ws.Cells.Item(892, 16).Value = 30000
$ws.Cells.Item(892, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(892, 18).Value = "Perú"
$ws.Cells.Item(892, 19).Value = 3000
$ws.Cells.Item(892, 20).Value = 10

# --- Row 893: Segunda ---
$ws.Cells.Item(893, 1).Value = 11
$ws.Cells.Item(893, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(893, 3).Value = "Bíobío"
$ws.Cells.Item(893, 4).Value = 45041
$ws.Cells.Item(893, 5).Value = 8
$ws.Cells.Item(893, 6).Value = "Fruta"
$ws.Cells.Item(893, 7).Value = 100106
$ws.Cells.Item(893, 8).Value = "Oleaginosos"
$ws.Cells.Item(893, 9).Value = 100106002
$ws.Cells.Item(893, 10).Value = "Palta"
$ws.Cells.Item(893, 11).Value = "Hass"
$ws.Cells.Item(893, 12).Value = "Segunda"
$ws.Cells.Item(893, 13).Value = 50
$ws.Cells.Item(893, 14).Value = 25000
$ws.Cells.Item(893, 15).Value = 25000
$ws.Cells.Item(893, 16).Value = 25000
$ws.Cells.Item(893, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(893, 18).Value = "Perú"
$ws.Cells.Item(893, 19).Value = 2500
$ws.Cells.Item(893, 20).Value = 10
